# Update res_bus vm_pu results for the 380 kV case: slack voltage setpoint
# dropped from 1.05 p.u. to 1.02 p.u., with the dependent bus voltages
# (columns B-F, I-N, rows 2-25) recomputed accordingly. Columns G/H and
# column A (bus indices) are unaffected by the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

# Row r => new values for columns B,C,D,E,F,I,J,K,L,M,N (in that order)
$newValues = @{
    2 = @(1.02, 1.024040460224237, 1.028496038485717, 1.027665297947351, 1.022515992300518, 1.031663216818477, 1.029217301905035, 1.031312686883958, 1.030484364629102, 1.025350143973158, 1.01377941020637)
    3 = @(1.02, 1.024889188612673, 1.029105872936048, 1.028461077837401, 1.023999682693419, 1.031832701924532, 1.029705288333132, 1.031731003574622, 1.031087949892214, 1.026638668804998, 1.013942538246503)
    4 = @(1.02, 1.025438607123161, 1.029500650977918, 1.028976584168031, 1.024959991830948, 1.031941284243814, 1.030020646449017, 1.032001163746147, 1.031478441090968, 1.027472199968253, 1.014047915641866)
    5 = @(1.02, 1.025669638111786, 1.029666656345226, 1.029193441766616, 1.025363772377986, 1.031986671962122, 1.030153126364585, 1.032114614317151, 1.031642586454573, 1.027822564441884, 1.014092173633045)
    6 = @(1.02, 1.025708432462715, 1.029694531731154, 1.029229861219261, 1.025431572923764, 1.031994277486775, 1.03017536464513, 1.032133655833828, 1.031670146159313, 1.027881389219931, 1.014099602236534)
    7 = @(1.02, 1.025441693953535, 1.029502868991272, 1.028979481285946, 1.024965386899995, 1.031941891739892, 1.030022417032429, 1.032002680169212, 1.031480634477516, 1.027476881755978, 1.014048507187219)
    8 = @(1.02, 1.02432724298922, 1.028702098148131, 1.027934114086769, 1.023017360406094, 1.031720719744845, 1.029382302108687, 1.031454166015846, 1.030688362950532, 1.025785656620624, 1.0138345766898)
    9 = @(1.02, 1.022365264088429, 1.027292419921237, 1.026096554757183, 1.01958652712302, 1.031322684269613, 1.02825128105463, 1.030483667066894, 1.029291780942482, 1.02280359414272, 1.013456254513323)
    10 = @(1.02, 1.021058544724684, 1.026353625939112, 1.024874603618385, 1.017300301730312, 1.03105176414241, 1.027495237457536, 1.029834047658989, 1.028360428610037, 1.020814070126477, 1.013203144247342)
    11 = @(1.02, 1.020493028805209, 1.025947365735778, 1.024346229020104, 1.016310523700627, 1.030933136978069, 1.027167385809665, 1.029552140565919, 1.0279570796008, 1.01995218777388, 1.013093334637871)
    12 = @(1.02, 1.020283016952952, 1.025796500231038, 1.024150078946634, 1.015942897323536, 1.030888875972889, 1.027045535371746, 1.0294473354292, 1.027807248253549, 1.019631981456523, 1.013052514899199)
    13 = @(1.02, 1.020328063120648, 1.025828859675437, 1.024192148731355, 1.016021753542974, 1.030898379046076, 1.027071675945139, 1.02946982065518, 1.027839388010319, 1.019700669795911, 1.013061272297187)
    14 = @(1.02, 1.020475668219715, 1.025934894372839, 1.024330012908668, 1.016280135183203, 1.030929482378913, 1.027157315075919, 1.029543479220456, 1.027944694686901, 1.019925720765264, 1.013089961108389)
    15 = @(1.02, 1.020566618696417, 1.026000230864704, 1.024414970415614, 1.016439335315289, 1.030948619993527, 1.027210070675041, 1.029588850479007, 1.028009576354552, 1.020064373385631, 1.013107633054789)
    16 = @(1.02, 1.021096082566622, 1.026380593312392, 1.024909685728778, 1.017365993276231, 1.0310596093123, 1.027516985812593, 1.029852743934666, 1.028387196184023, 1.020871261543775, 1.013210427514322)
    17 = @(1.02, 1.021428282549998, 1.026619250601342, 1.02522020584638, 1.017947304253452, 1.031128877578218, 1.027709377362202, 1.03001811232985, 1.028624049493804, 1.021377289872452, 1.0132748513402)
    18 = @(1.02, 1.021622078300626, 1.026758478815495, 1.025401398156351, 1.018286390084751, 1.031169153500418, 1.027821549783296, 1.030114509294075, 1.028762195492002, 1.021672408792049, 1.013312408279452)
    19 = @(1.02, 1.021688162566389, 1.026805955986524, 1.025463192098696, 1.018402012611818, 1.03118286498514, 1.027859789813808, 1.0301473680395, 1.028809298583094, 1.021773030248354, 1.013325210761784)
    20 = @(1.02, 1.021392637619626, 1.026593642499528, 1.025186882607876, 1.017884933364783, 1.031121458891938, 1.027688740357718, 1.03000037602951, 1.028598638046663, 1.021323001867547, 1.013267941381012)
    21 = @(1.02, 1.020432200924511, 1.025903668750727, 1.02428941225215, 1.016204047681149, 1.030920328675077, 1.027132098463989, 1.029521791152563, 1.027913684755056, 1.019859450685528, 1.013081513838412)
    22 = @(1.02, 1.01982860231423, 1.025470073037195, 1.023725784153315, 1.015147328772917, 1.030792727034875, 1.026781700724465, 1.029220352125351, 1.027482972355656, 1.018938880888819, 1.012964116926751)
    23 = @(1.02, 1.020148555958969, 1.025699909317688, 1.024024512447432, 1.015707505593236, 1.030860479342242, 1.026967492408666, 1.029380201082362, 1.027711306187689, 1.019426929417347, 1.013026368495222)
    24 = @(1.02, 1.021408743939346, 1.026605213626057, 1.025201939725085, 1.017913116031089, 1.031124811469608, 1.027698065474815, 1.030008390484418, 1.028610120405981, 1.021347532399221, 1.01327106375713)
    25 = @(1.02, 1.022872262256203, 1.027656685300001, 1.026571067568537, 1.020473286520871, 1.031426567896093, 1.028544036698059, 1.030735028567717, 1.029652886057318, 1.023574775393616, 1.013554218551669)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $vals[$i]
    }
}
